$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4729.1763
$ws.Range("J2").Value = 6042.143
$ws.Range("L2").Value = 6042.143
$ws.Range("N2").Value = -6268.143
$ws.Range("H15").Value = 1022.381
$ws.Range("I15").Value = 1022.381
$ws.Range("K15").Value = 3067.143
$ws.Range("M15").Value = -2898.143
$ws.Range("H38").Value = 1428.3334
$ws.Range("I38").Value = 24.666666
$ws.Range("J38").Value = 2832
$ws.Range("K38").Value = 73.99999800000001
$ws.Range("L38").Value = 8496
$ws.Range("M38").Value = 298.000002
$ws.Range("N38").Value = -9240
$ws.Range("H55").Value = 1164.9166
$ws.Range("I55").Value = 247.5
$ws.Range("K55").Value = 247.5
$ws.Range("M55").Value = -33.5
$ws.Range("H86").Value = 102566456
$ws.Range("J86").Value = 111112070
$ws.Range("L86").Value = 111112070
$ws.Range("N86").Value = -111114316
$ws.Range("H88").Value = 7353.2
$ws.Range("J88").Value = 7353.2
$ws.Range("L88").Value = 7353.2
$ws.Range("N88").Value = -8165.2
$ws.Range("H89").Value = 102566456
$ws.Range("J89").Value = 111112070
$ws.Range("L89").Value = 555560350
$ws.Range("N89").Value = -555571582
$ws.Range("H91").Value = 7353.2
$ws.Range("J91").Value = 7353.2
$ws.Range("L91").Value = 7353.2
$ws.Range("N91").Value = -10161.2
$ws.Range("H92").Value = 1674555.2
$ws.Range("I92").Value = 710738.1
$ws.Range("K92").Value = 710738.1
$ws.Range("M92").Value = -709490.1
$ws.Range("H96").Value = 1562.3636
$ws.Range("I96").Value = 751.4
$ws.Range("J96").Value = 2238.1667
$ws.Range("K96").Value = 2254.2
$ws.Range("L96").Value = 6714.500100000001
$ws.Range("M96").Value = -881.1999999999998
$ws.Range("N96").Value = -9460.500100000001
$ws.Range("H98").Value = 3841.9443
$ws.Range("I98").Value = 2747.5
$ws.Range("K98").Value = 2747.5
$ws.Range("M98").Value = -1249.5
$ws.Range("H99").Value = 813.3333
$ws.Range("J99").Value = 1100
$ws.Range("L99").Value = 3300
$ws.Range("N99").Value = -6296
$ws.Range("H103").Value = 649.8333
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 719.8
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 2159.4
$ws.Range("M103").Value = -314
$ws.Range("N103").Value = -3331.4
$ws.Range("H115").Value = 731.4286
$ws.Range("I115").Value = 694
$ws.Range("J115").Value = 825
$ws.Range("K115").Value = 2082
$ws.Range("L115").Value = 2475
$ws.Range("M115").Value = -515
$ws.Range("N115").Value = -5609
$ws.Range("H122").Value = 3841.9443
$ws.Range("I122").Value = 2747.5
$ws.Range("K122").Value = 8242.5
$ws.Range("M122").Value = -5792.5
$ws.Range("H129").Value = 2102.125
$ws.Range("I129").Value = 1500
$ws.Range("K129").Value = 4500
$ws.Range("M129").Value = 500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 653.63635
$ws.Range("I2").Value = 635
$ws.Range("J2").Value = 680.55554
$ws.Range("K2").Value = 635
$ws.Range("L2").Value = 680.55554
$ws.Range("M2").Value = -522
$ws.Range("N2").Value = -906.55554
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 9996.625
$ws.Range("I36").Value = 9996
$ws.Range("K36").Value = 9996
$ws.Range("M36").Value = -9650
$ws.Range("H61").Value = 3582253.5
$ws.Range("I61").Value = 82337.53999999999
$ws.Range("K61").Value = 82337.53999999999
$ws.Range("M61").Value = -82125.53999999999
$ws.Range("H102").Value = 2002.6842
$ws.Range("I102").Value = 2002.6842
$ws.Range("K102").Value = 2002.6842
$ws.Range("M102").Value = -380.6841999999999
$ws.Range("H116").Value = 653.63635
$ws.Range("I116").Value = 635
$ws.Range("J116").Value = 680.55554
$ws.Range("K116").Value = 635
$ws.Range("L116").Value = 680.55554
$ws.Range("M116").Value = 1659
$ws.Range("N116").Value = -5268.55554
$ws.Range("H122").Value = 1510.3158
$ws.Range("I122").Value = 1452.7059
$ws.Range("K122").Value = 4358.1177
$ws.Range("M122").Value = -1908.1177
$ws.Range("H132").Value = 2317.6333
$ws.Range("I132").Value = 2075.4348
$ws.Range("K132").Value = 6226.3044
$ws.Range("M132").Value = -3696.3044
$ws.Range("H136").Value = 3582253.5
$ws.Range("I136").Value = 82337.53999999999
$ws.Range("K136").Value = 247012.62
$ws.Range("M136").Value = -244462.62

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 653.63635
$ws.Range("I3").Value = 635
$ws.Range("J3").Value = 680.55554
$ws.Range("K3").Value = 635
$ws.Range("L3").Value = 680.55554
$ws.Range("M3").Value = -521
$ws.Range("N3").Value = -908.55554
$ws.Range("H100").Value = 37799.8
$ws.Range("J100").Value = 37799.8
$ws.Range("L100").Value = 37799.8
$ws.Range("N100").Value = -39963.8
$ws.Range("H134").Value = 21953028
$ws.Range("I134").Value = 1548.3715
$ws.Range("K134").Value = 4645.1145
$ws.Range("M134").Value = -2110.1145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2882.0732
$ws.Range("I31").Value = 4085.5
$ws.Range("K31").Value = 4085.5
$ws.Range("M31").Value = -3790.5
$ws.Range("H34").Value = 2882.0732
$ws.Range("I34").Value = 4085.5
$ws.Range("K34").Value = 4085.5
$ws.Range("M34").Value = -3883.5
$ws.Range("H99").Value = 86668990
$ws.Range("I99").Value = 6668994.5
$ws.Range("J99").Value = 166669000
$ws.Range("K99").Value = 6668994.5
$ws.Range("L99").Value = 166669000
$ws.Range("M99").Value = -6667496.5
$ws.Range("N99").Value = -166671996
$ws.Range("H122").Value = 4972.0625
$ws.Range("I122").Value = 4782.5
$ws.Range("K122").Value = 14347.5
$ws.Range("M122").Value = -11897.5
$ws.Range("H126").Value = 86668990
$ws.Range("I126").Value = 6668994.5
$ws.Range("J126").Value = 166669000
$ws.Range("K126").Value = 20006983.5
$ws.Range("L126").Value = 500007000
$ws.Range("M126").Value = -20004513.5
$ws.Range("N126").Value = -500011940
$ws.Range("H132").Value = 10787309
$ws.Range("I132").Value = 41949.2
$ws.Range("J132").Value = 55559644
$ws.Range("K132").Value = 125847.6
$ws.Range("L132").Value = 166678932
$ws.Range("M132").Value = -123317.6
$ws.Range("N132").Value = -166683992
$ws.Range("H134").Value = 2374.2917
$ws.Range("I134").Value = 1984.9375
$ws.Range("K134").Value = 5954.8125
$ws.Range("M134").Value = -3419.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 200310.6
$ws.Range("I23").Value = 351
$ws.Range("J23").Value = 500250
$ws.Range("K23").Value = 1053
$ws.Range("L23").Value = 1500750
$ws.Range("M23").Value = -818
$ws.Range("N23").Value = -1501220
$ws.Range("H138").Value = 6034
$ws.Range("I138").Value = 5875
$ws.Range("J138").Value = 6749.5
$ws.Range("K138").Value = 17625
$ws.Range("L138").Value = 20248.5
$ws.Range("M138").Value = -12485
$ws.Range("N138").Value = -30528.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 5866318
$ws.Range("I132").Value = 1649.1052
$ws.Range("K132").Value = 4947.3156
$ws.Range("M132").Value = -2417.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3288.75
$ws.Range("I61").Value = 3177.6667
$ws.Range("K61").Value = 3177.6667
$ws.Range("M61").Value = -2975.6667
$ws.Range("H68").Value = 4949.5
$ws.Range("J68").Value = 4949.5
$ws.Range("L68").Value = 4949.5
$ws.Range("N68").Value = -6447.5
$ws.Range("H71").Value = 4949.5
$ws.Range("J71").Value = 4949.5
$ws.Range("L71").Value = 24747.5
$ws.Range("N71").Value = -32235.5
$ws.Range("H93").Value = 1372.909
$ws.Range("J93").Value = 2450
$ws.Range("L93").Value = 2450
$ws.Range("N93").Value = -4946
$ws.Range("H100").Value = 3339.3125
$ws.Range("I100").Value = 2817.875
$ws.Range("J100").Value = 3860.75
$ws.Range("K100").Value = 2817.875
$ws.Range("L100").Value = 3860.75
$ws.Range("M100").Value = -2276.875
$ws.Range("N100").Value = -4942.75
$ws.Range("H113").Value = 3288.75
$ws.Range("I113").Value = 3177.6667
$ws.Range("K113").Value = 3177.6667
$ws.Range("M113").Value = -1007.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1436
$ws.Range("I100").Value = 1195
$ws.Range("K100").Value = 2390
$ws.Range("M100").Value = -1849
$ws.Range("H122").Value = 1862.0769
$ws.Range("I122").Value = 1862.0769
$ws.Range("K122").Value = 5586.2307
$ws.Range("M122").Value = -3136.2307
$ws.Range("H132").Value = 2363.2222
$ws.Range("I132").Value = 2258.2778
$ws.Range("K132").Value = 6774.8334
$ws.Range("M132").Value = -4244.8334
